# Generate Report for Handoff
#
# The localization run finished: status flips from "In Translation" to
# "Ready for handoff", and the HO-xliff-generate / handoff timestamps move
# forward a few seconds. Excel then re-autofits the "Status"/"zh-cn"/"de-de"
# columns because "Ready for handoff" is wider than "In Translation".

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Overview sheet (row 2): zh-cn / de-de status columns + HO xliff date ---
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Range("G2").Value = "2016-08-28 22:39:37"

# --- zh-cn detail sheet (row 2): Status + Latest Handoff Datetime ---
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("H2").Value = "2016-08-28 22:39:33"

# --- de-de detail sheet (row 2): Status + Latest Handoff Datetime ---
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("H2").Value = "2016-08-28 22:39:37"

# --- Column widths auto-adjust to fit the new, longer "Ready for handoff" text ---
$wsOverview.Columns.Item(5).ColumnWidth = 16.33   # E: zh-cn
$wsOverview.Columns.Item(6).ColumnWidth = 16.33   # F: de-de
$wsZhCn.Columns.Item(3).ColumnWidth = 16.33        # C: Status
$wsDeDe.Columns.Item(3).ColumnWidth = 16.33        # C: Status
